# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late" / "heading" / "Outstanding")
#   columns one place to the right (-> O/P/Q).
# - Give the newly inserted column the same width the author set by hand.
# - Make "Repayment schedule" the active sheet and leave the selection on N15,
#   which is what moves tabSelected/activeTab off "NewLoanInput" and onto
#   "Repayment schedule".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column at N; this shifts existing N:P -> O:Q and keeps every
# cell's style/value intact.
$ws.Columns("N").Insert()

# The author widened the newly-inserted column (it's no longer an
# autofit/bestFit column like its neighbours).
$ws.Columns("N").ColumnWidth = 10.2

# Switch focus to the "Repayment schedule" tab and park the selection on N15.
$ws.Activate() | Out-Null
$ws.Range("N15").Select() | Out-Null
